$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New reading link for row 7
$ws.Range("A7").Value = "[Models of network structure](#sec:models)"

# 2. Move the "Small worlds" reading link from A9 down to A10
$ws.Range("A9").Value = $null
$ws.Range("A10").Value = "[Small worlds](#sec:smallworlds)"

# 3. New reading link for row 20
$ws.Range("A20").Value = "[Concurrency](#sec:concurrency)"

# 4. Rename the reading link in row 22
$ws.Range("A22").Value = "[Social influence](#sec:socialinfluence)"

# 5. The old row-22 reading link now appears at row 24
$ws.Range("A24").Value = "[Dynamics: Complex contagion and social influence](#sec:complexcontagion)"
